$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New shuttlecock purchase record (row 16 / index 15), added 2024-03-19
$ws.Range("A16").Value = 15

# Copy the date-formatted style from the row above, then set the value
$ws.Range("B15").Copy($ws.Range("B16"))
$ws.Range("B16").Value = 45370

$ws.Range("C16").Value = "Bullet tournament 76"
$ws.Range("D16").Value = 670
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 99

$ws.Range("G16").Formula = "=D16*E16+F16"
$ws.Range("H16").Formula = "=E16*12"
$ws.Range("I16").Formula = "=ROUNDUP(G16/H16,0)"

# Update the selected cell to reflect where the user left off editing
$ws.Range("E32").Select()
